# 2.1.1.1e.xlsx update: add a 2020 column (N) and refresh several
# previously-estimated 2018/2019 figures (columns L/M) with final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting into the new column N so it matches the look of
#     column M / the existing thick-bordered rows (row 3 separator, row 4
#     header, rows 5-13 data, row 14 footer). PasteSpecial(-4122) ==
#     xlPasteFormats, so only styles move, not cell contents.
$ws.Range("A14").Copy()
$ws.Range("N3").PasteSpecial(-4122)

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)

$ws.Range("M5:M13").Copy()
$ws.Range("N5:N13").PasteSpecial(-4122)

$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- New 2020 header value
$ws.Range("N4").Value = 2020

# --- Updated 2018 figures
$ws.Range("L9").Value = 105.6
$ws.Range("L12").Value = 27.1

# --- Updated 2019 figures (column M)
$ws.Range("M5").Value = 68.400000000000006
$ws.Range("M6").Value = 108.2
$ws.Range("M7").Value = 51.7
$ws.Range("M8").Value = 97.7
$ws.Range("M9").Value = 106.7
$ws.Range("M10").Value = 124.2
$ws.Range("M11").Value = 138.80000000000001
$ws.Range("M12").Value = 33.9
$ws.Range("M13").Value = 96
$ws.Range("M14").Value = 7.7

# --- New 2020 figures (column N)
$ws.Range("N5").Value = 68.5
$ws.Range("N6").Value = 106.7
$ws.Range("N7").Value = 53.2
$ws.Range("N8").Value = 49.6
$ws.Range("N9").Value = 108.9
$ws.Range("N10").Value = 107.8
$ws.Range("N11").Value = 155.69999999999999
$ws.Range("N12").Value = 25.9
$ws.Range("N13").Value = 103.5
$ws.Range("N14").Value = 11

# --- Page setup (Page Setup > Paper Size = A4 / Letter #9, Portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
